$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns for each data row, in order: A, D, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
# (B is always "Agrn", C is always "Atp1a3")
$rows = @(
  @("ECs", "ECs", 3, 1, 17.306265, 51.918795, 0.5463168539988408, 0.5463168539988407, 1, 0.3333333333333333, 11.76781466666667, 35.303444, 0.9616440602480689, 0.9616440602480689, 203.65691909222, 1832.91227182998, 0.5253623576613967, 0.5253623576613966),
  @("ECs", "FAPs", 3, 1, 17.306265, 51.918795, 0.5463168539988408, 0.5463168539988407, 3, 1, 0.3773413333333333, 1.132024, 0.03083563619623796, 0.03083563619623796, 6.53036911012, 58.77332199108, 0.0168460277577815, 0.0168460277577815),
  @("ECs", "sCs", 3, 1, 17.306265, 51.918795, 0.5463168539988408, 0.5463168539988407, 2, 0.6666666666666666, 0.09202733333333334, 0.276082, 0.007520303555692963, 0.007520303555692962, 1.59264941791, 14.33384476119, 0.004108468579662476, 0.004108468579662474),
  @("FAPs", "ECs", 3, 1, 3.344413333333334, 10.03324, 0.1055750256186672, 0.1055750256186672, 1, 0.3333333333333333, 11.76781466666667, 35.303444, 0.9616440602480689, 0.9616440602480689, 39.35643627539556, 354.20792647856, 0.101525596296729, 0.101525596296729),
  @("FAPs", "FAPs", 3, 1, 3.344413333333334, 10.03324, 0.1055750256186672, 0.1055750256186672, 3, 1, 0.3773413333333333, 1.132024, 0.03083563619623796, 0.03083563619623796, 1.261985386417778, 11.35786847776, 0.003255473081385724, 0.003255473081385723),
  @("FAPs", "sCs", 3, 1, 3.344413333333334, 10.03324, 0.1055750256186672, 0.1055750256186672, 2, 0.6666666666666666, 0.09202733333333334, 0.276082, 0.007520303555692963, 0.007520303555692962, 0.3077774406311112, 2.76999696568, 0.0007939562405524384, 0.0007939562405524382),
  @("sCs", "ECs", 3, 1, 11.02739433333333, 33.082183, 0.3481081203824922, 0.3481081203824921, 1, 0.3333333333333333, 11.76781466666667, 35.303444, 0.9616440602480689, 0.9616440602480689, 129.7683327709169, 1167.914994938252, 0.3347561062899433, 0.3347561062899433),
  @("sCs", "FAPs", 3, 1, 11.02739433333333, 33.082183, 0.3481081203824922, 0.3481081203824921, 3, 1, 0.3773413333333333, 1.132024, 0.03083563619623796, 0.03083563619623796, 4.161091680932445, 37.449825128392, 0.01073413535707074, 0.01073413535707074),
  @("sCs", "sCs", 3, 1, 11.02739433333333, 33.082183, 0.3481081203824922, 0.3481081203824921, 2, 0.6666666666666666, 0.09202733333333334, 0.276082, 0.007520303555692963, 0.007520303555692962, 1.014821694111778, 9.133395247006, 0.00261787873547805, 0.002617878735478049)
)

$rowIndex = 2
foreach ($r in $rows) {
  $ws.Cells.Item($rowIndex, 1).Value = $r[0]          # A: Sending cluster
  $ws.Cells.Item($rowIndex, 2).Value = "Agrn"          # B: Ligand symbol
  $ws.Cells.Item($rowIndex, 3).Value = "Atp1a3"        # C: Receptor symbol
  $ws.Cells.Item($rowIndex, 4).Value = $r[1]          # D: Target cluster
  $ws.Cells.Item($rowIndex, 5).Value = $r[2]          # E
  $ws.Cells.Item($rowIndex, 6).Value = $r[3]          # F
  $ws.Cells.Item($rowIndex, 7).Value = $r[4]          # G
  $ws.Cells.Item($rowIndex, 8).Value = $r[5]          # H
  $ws.Cells.Item($rowIndex, 9).Value = $r[6]          # I
  $ws.Cells.Item($rowIndex, 10).Value = $r[7]         # J
  $ws.Cells.Item($rowIndex, 11).Value = $r[8]         # K
  $ws.Cells.Item($rowIndex, 12).Value = $r[9]         # L
  $ws.Cells.Item($rowIndex, 13).Value = $r[10]        # M
  $ws.Cells.Item($rowIndex, 14).Value = $r[11]        # N
  $ws.Cells.Item($rowIndex, 15).Value = $r[12]        # O
  $ws.Cells.Item($rowIndex, 16).Value = $r[13]        # P
  $ws.Cells.Item($rowIndex, 17).Value = $r[14]        # Q
  $ws.Cells.Item($rowIndex, 18).Value = $r[15]        # R
  $ws.Cells.Item($rowIndex, 19).Value = $r[16]        # S
  $ws.Cells.Item($rowIndex, 20).Value = $r[17]        # T
  $rowIndex++
}
